$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.964.97"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.889.32"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  +1.30%  "
$ws.Range("D5").Value = "'336.18"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'1.014"
$ws.Range("E6").Value = "  +1.21%  "
$ws.Range("D7").Value = "'0.4722"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'0.3946"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "'47.03"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "'0.08013"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'1.019"
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "'21.81"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "1.901.06"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'7.190"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'0.06795"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("D18").Value = "'88.14"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").Value = "'0.00001053"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "'17.17"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D22").Value = "27.950.57"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "'5.508"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'11.00"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("D25").Value = "'2.361"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "2.118.77"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'159.45"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "'20.04"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "'2.110"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "'5.509"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").Value = "'121.73"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").Value = "'0.09578"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "'0.9670"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").Value = "'3.649"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").Value = "'5.364"
$ws.Range("E35").Value = "  +0.46%  "
$ws.Range("D36").Value = "'1.363"
$ws.Range("E36").Value = "  -7.15%  "
$ws.Range("D37").Value = "'0.06123"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").Value = "'8.305"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'0.5961"
$ws.Range("E41").Value = "  -1.23%  "
$ws.Range("D42").Value = "'0.1910"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "'10.41"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "'1.276"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.5704"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'1.948"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'3.399"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").Value = "'0.06874"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'113.47"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'1.071"
$ws.Range("E51").Value = "  -0.45%  "
